$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task-progress rows appended below the existing data (rows 2-12),
# matching the date/status formatting already used in that block.
$rows = @(
    @(45610, "Created Config.py", "Complete", "Adriean"),
    @(45610, "Created Individuals.py", "Complete", "Adriean"),
    @(45615, "Set up dictionary for configs in Individuals.py", "Complete", "Adriean"),
    @(45615, "Worked on csv import methods", "Complete", "Nicholas"),
    @(45615, "Worked on dataframe objects in Individuals.py", "Complete", "Nicholas"),
    @(45615, "Worked on line graphs in Individuals.py", "Complete", "Nicholas"),
    @(45616, "Worked on query functionality in Individuals.py", "Complete", "Chris"),
    @(45616, "Worked on violin plots in Individuals.py", "Complete", "Chris"),
    @(45616, "Cleaned up Individuals.py", "Complete", "Chris"),
    @(45616, "Tested method functionality of Individuals.py", "Complete", "Chris")
)

$startRow = 13
$endRow = $startRow + $rows.Length - 1

# Carry the date-column number format down from the last existing data
# row (A12) before writing the new values into A13:A22.
$ws.Range("A12").Copy()
$ws.Range("A13:A" + $endRow).PasteSpecial(-4122)

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("B20").Select()
